$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.998.69"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -2.72%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.620.65"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.10%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "512.54"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.29"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -4.07%  "
$ws.Range("E7").Value = "  -0.40%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.576"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.70%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.645.77"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.50%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.47"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.106"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.339"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.67%  "
$ws.Range("E13").Value = "  -0.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.076.25"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "58.767.26"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -3.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.30"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.97%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000139"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.95%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.627.12"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.61"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "346.97"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.51"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.21"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.95"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.425"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.724.69"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.993"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.59%  "
$ws.Range("E28").Value = "  -3.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0833"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.11"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.997"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.35"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.12"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.56"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "149.52"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.994"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +12.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.05"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.15"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.880"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.32"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.43%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.44"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.66"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "291.62"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -4.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.623"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.63%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0999"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.993"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.74"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0543"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.78"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0232"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "10.25"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.99%  "
